# Update the author byline from "Adam Zabell" to "Randall Julian".
# The paragraph (style "Author") reads "Adam Zabell"; swap each name in
# place with a whole-word Find/Replace so the surrounding space and the
# rest of the document are left untouched.

$d = $word.ActiveDocument

$d.Content.Find.Execute("Adam", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Randall", 2) | Out-Null

$d.Content.Find.Execute("Zabell", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Julian", 2) | Out-Null
